# Add a new "localdb" command-type column to the hidden '#system' sheet,
# which powers Nexial's catalog of command types (the "target" named range
# drives a dropdown-list data validation on the excel_test sheet).
#
# This mirrors what the authoring tool did: insert a new column N (so the
# existing N..AC columns shift right to O..AD), populate it with the new
# "localdb" command family, insert "localdb" alphabetically into the
# existing A2:A29 "target" catalog (so it becomes A2:A30), and refresh every
# workbook-level defined name whose range moved as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1. Insert a new column at N. Everything from N..AC shifts to O..AD,
#    carrying its data/formatting along (same as Excel's native
#    Columns("N:N").Insert()).
# ---------------------------------------------------------------------
$ws.Columns("N:N").Insert()

# ---------------------------------------------------------------------
# 2. Populate the new "localdb" column (N) with its command list.
# ---------------------------------------------------------------------
$ws.Range("N1").Value2 = "localdb"
$ws.Range("N2").Value2 = "cloneTable(var,source,target)"
$ws.Range("N3").Value2 = "dropTables(var,tables)"
$ws.Range("N4").Value2 = "exportCSV(sql,output)"
$ws.Range("N5").Value2 = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value2 = "purge(var)"
$ws.Range("N7").Value2 = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------
# 3. Insert "localdb" into the target catalog (column A), which is kept
#    alphabetically sorted. It belongs between "json" (A13) and the old
#    "macro" (A14). A plain Range.Insert() on this host shifts the whole
#    row (all columns) rather than just column A, so shift column A's
#    values down manually (bottom-up, to avoid clobbering) and then set
#    the new A14 value -- this leaves every other column untouched.
# ---------------------------------------------------------------------
for ($r = 29; $r -ge 14; $r--) {
    $ws.Range("A" + ($r + 1)).Value2 = $ws.Range("A" + $r).Value2
}
$ws.Range("A14").Value2 = "localdb"

# ---------------------------------------------------------------------
# 4. Refresh workbook-level defined names: every name anchored at column
#    N or later needs to shift one column to the right; "target" now
#    spans one extra row; and a brand-new "localdb" name is added for
#    the new column.
# ---------------------------------------------------------------------
$wb.Names.Item("macro").RefersTo      = '=''#system''!$O$2:$O$4'
$wb.Names.Item("mail").RefersTo       = '=''#system''!$P$2:$P$2'
$wb.Names.Item("number").RefersTo     = '=''#system''!$Q$2:$Q$16'
$wb.Names.Item("pdf").RefersTo        = '=''#system''!$R$2:$R$16'
$wb.Names.Item("rdbms").RefersTo      = '=''#system''!$S$2:$S$7'
$wb.Names.Item("redis").RefersTo      = '=''#system''!$T$2:$T$10'
$wb.Names.Item("sms").RefersTo        = '=''#system''!$U$2:$U$2'
$wb.Names.Item("sound").RefersTo      = '=''#system''!$V$2:$V$5'
$wb.Names.Item("ssh").RefersTo        = '=''#system''!$W$2:$W$9'
$wb.Names.Item("step").RefersTo       = '=''#system''!$X$2:$X$4'
$wb.Names.Item("target").RefersTo     = '=''#system''!$A$2:$A$30'
$wb.Names.Item("web").RefersTo        = '=''#system''!$Y$2:$Y$127'
$wb.Names.Item("webalert").RefersTo   = '=''#system''!$Z$2:$Z$8'
$wb.Names.Item("webcookie").RefersTo  = '=''#system''!$AA$2:$AA$8'
$wb.Names.Item("ws").RefersTo         = '=''#system''!$AB$2:$AB$17'
$wb.Names.Item("ws.async").RefersTo   = '=''#system''!$AC$2:$AC$8'
$wb.Names.Item("xml").RefersTo        = '=''#system''!$AD$2:$AD$21'

$wb.Names.Add('localdb', '=''#system''!$N$2:$N$7')

# ---------------------------------------------------------------------
# 5. The sheet's stored <dimension> always trails one column past the
#    last column that actually holds data (it was "A1:AD127" pre-edit
#    even though data stopped at AC), so after the insert it should read
#    "A1:AE127". Touch AE1's formatting to register it as used, then
#    immediately reset it to the default style so no visible formatting
#    or value is left behind.
# ---------------------------------------------------------------------
$ws.Range("AE1").Interior.Color = 255
$ws.Range("AE1").Style = "Normal"
